$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.956.56"
$ws.Range("E2").Value = "  +3.23%  "

$ws.Range("D3").Value = "1.912.06"
$ws.Range("E3").Value = "  +1.70%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.47%  "

$ws.Range("D5").Value = "'246.26"
$ws.Range("E5").Value = "  +1.00%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.57%  "

$ws.Range("D7").Value = "'0.4968"
$ws.Range("E7").Value = "  +0.52%  "

$ws.Range("D8").Value = "'0.2994"
$ws.Range("E8").Value = "  +2.90%  "

$ws.Range("D9").Value = "'0.06777"
$ws.Range("E9").Value = "  +2.51%  "

$ws.Range("D10").Value = "1.916.56"
$ws.Range("E10").Value = "  +1.93%  "

$ws.Range("D11").Value = "'17.02"
$ws.Range("E11").Value = "  +0.43%  "

$ws.Range("D12").Value = "'0.07335"
$ws.Range("E12").Value = "  +2.18%  "

$ws.Range("D13").Value = "'0.6840"
$ws.Range("E13").Value = "  +2.10%  "

$ws.Range("D14").Value = "'89.10"
$ws.Range("E14").Value = "  +3.99%  "

$ws.Range("D15").Value = "'5.072"
$ws.Range("E15").Value = "  +5.03%  "

$ws.Range("D16").Value = "30.894.30"
$ws.Range("E16").Value = "  +3.07%  "

$ws.Range("D17").Value = "'0.000008034"
$ws.Range("E17").Value = "  +2.12%  "

$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "'13.17"
$ws.Range("E19").Value = "  +3.10%  "

$ws.Range("D20").Value = "2.156.38"
$ws.Range("E20").Value = "  +1.87%  "

$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.64%  "

$ws.Range("D22").Value = "'4.875"
$ws.Range("E22").Value = "  +2.53%  "

$ws.Range("D23").Value = "'176.28"
$ws.Range("E23").Value = "  +30.97%  "

$ws.Range("D24").Value = "'6.044"
$ws.Range("E24").Value = "  +8.05%  "

$ws.Range("D25").Value = "'9.333"
$ws.Range("E25").Value = "  +2.41%  "

$ws.Range("D26").Value = "'152.48"
$ws.Range("E26").Value = "  +3.57%  "

$ws.Range("D27").Value = "'18.12"
$ws.Range("E27").Value = "  +8.34%  "

$ws.Range("D28").Value = "'1.949"
$ws.Range("E28").Value = "  +1.14%  "

$ws.Range("D29").Value = "'1.418"
$ws.Range("E29").Value = "  +2.96%  "

$ws.Range("D30").Value = "'4.327"
$ws.Range("E30").Value = "  +3.52%  "

$ws.Range("D31").Value = "'0.08905"
$ws.Range("E31").Value = "  +3.33%  "

$ws.Range("D32").Value = "'4.072"
$ws.Range("E32").Value = "  +3.49%  "

$ws.Range("D33").Value = "'0.05304"
$ws.Range("E33").Value = "  +6.12%  "

$ws.Range("D34").Value = "'0.7465"
$ws.Range("E34").Value = "  +5.48%  "

$ws.Range("D35").Value = "'1.140"
$ws.Range("E35").Value = "  +2.42%  "

$ws.Range("D36").Value = "'2.632"
$ws.Range("E36").Value = "  -0.57%  "

$ws.Range("D37").Value = "'0.01898"
$ws.Range("E37").Value = "  +15.72%  "

$ws.Range("D38").Value = "'2.721"
$ws.Range("E38").Value = "  +1.28%  "

$ws.Range("D39").Value = "'2.214"
$ws.Range("E39").Value = "  -0.27%  "

$ws.Range("D40").Value = "'0.9443"
$ws.Range("E40").Value = "  +1.24%  "

$ws.Range("D41").Value = "'5.994"

$ws.Range("D42").Value = "'0.4386"
$ws.Range("E42").Value = "  +4.74%  "

$ws.Range("D43").Value = "'105.34"
$ws.Range("E43").Value = "  +2.51%  "

$ws.Range("D44").Value = "'7.837"
$ws.Range("E44").Value = "  +2.33%  "

$ws.Range("E45").Value = "  +0.85%  "

$ws.Range("E46").Value = "  +5.33%  "

$ws.Range("D47").Value = "'0.05833"
$ws.Range("E47").Value = "  +2.29%  "

$ws.Range("D48").Value = "'0.3924"
$ws.Range("E48").Value = "  +5.74%  "

$ws.Range("D49").Value = "'33.41"
$ws.Range("E49").Value = "  +2.32%  "

$ws.Range("D50").Value = "'8.526"
$ws.Range("E50").Value = "  +4.20%  "

$ws.Range("D51").Value = "'1.381"
$ws.Range("E51").Value = "  +3.19%  "
